{"js": "// Replace each \"three-digit number \u00d7 one-digit number =\" expression in the\n// worksheet table with a newly generated problem, matching the commit's\n// regenerated numbers. Every old expression is unique in the document, so a\n// scoped search-and-replace (old text -> new text) for each pair is\n// unambiguous and leaves everything else (fonts, sizes, paragraph layout,\n// the date header, blank drill rows) untouched.\nconst replacements = [\n  [\"305\u00d73=\", \"292\u00d75=\"],\n  [\"470\u00d79=\", \"757\u00d75=\"],\n  [\"957\u00d77=\", \"602\u00d74=\"],\n  [\"296\u00d73=\", \"220\u00d78=\"],\n  [\"770\u00d78=\", \"549\u00d72=\"],\n  [\"324\u00d79=\", \"883\u00d75=\"],\n  [\"628\u00d75=\", \"929\u00d74=\"],\n  [\"231\u00d76=\", \"478\u00d74=\"],\n  [\"814\u00d78=\", \"625\u00d74=\"],\n  [\"780\u00d75=\", \"734\u00d74=\"],\n  [\"310\u00d72=\", \"941\u00d73=\"],\n  [\"817\u00d79=\", \"857\u00d76=\"],\n  [\"180\u00d74=\", \"927\u00d79=\"],\n  [\"965\u00d78=\", \"937\u00d78=\"],\n  [\"351\u00d78=\", \"149\u00d73=\"],\n  [\"443\u00d77=\", \"829\u00d77=\"],\n  [\"724\u00d78=\", \"554\u00d78=\"],\n  [\"542\u00d75=\", \"832\u00d77=\"],\n  [\"799\u00d74=\", \"430\u00d73=\"],\n  [\"525\u00d78=\", \"134\u00d76=\"],\n  [\"758\u00d76=\", \"445\u00d76=\"],\n  [\"622\u00d78=\", \"595\u00d77=\"],\n  [\"349\u00d75=\", \"539\u00d73=\"],\n  [\"488\u00d79=\", \"396\u00d73=\"],\n  [\"444\u00d76=\", \"308\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"three-digit number x one-digit number =\" expression in the\n# worksheet table with a newly generated problem, matching the commit's\n# regenerated numbers. Every old expression is unique in the document, so a\n# simple Find/Replace per pair (old text -> new text) is unambiguous and\n# leaves everything else (fonts, sizes, paragraph layout, the date header,\n# blank drill rows) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"305\u00d73=\", \"292\u00d75=\"),\n    @(\"470\u00d79=\", \"757\u00d75=\"),\n    @(\"957\u00d77=\", \"602\u00d74=\"),\n    @(\"296\u00d73=\", \"220\u00d78=\"),\n    @(\"770\u00d78=\", \"549\u00d72=\"),\n    @(\"324\u00d79=\", \"883\u00d75=\"),\n    @(\"628\u00d75=\", \"929\u00d74=\"),\n    @(\"231\u00d76=\", \"478\u00d74=\"),\n    @(\"814\u00d78=\", \"625\u00d74=\"),\n    @(\"780\u00d75=\", \"734\u00d74=\"),\n    @(\"310\u00d72=\", \"941\u00d73=\"),\n    @(\"817\u00d79=\", \"857\u00d76=\"),\n    @(\"180\u00d74=\", \"927\u00d79=\"),\n    @(\"965\u00d78=\", \"937\u00d78=\"),\n    @(\"351\u00d78=\", \"149\u00d73=\"),\n    @(\"443\u00d77=\", \"829\u00d77=\"),\n    @(\"724\u00d78=\", \"554\u00d78=\"),\n    @(\"542\u00d75=\", \"832\u00d77=\"),\n    @(\"799\u00d74=\", \"430\u00d73=\"),\n    @(\"525\u00d78=\", \"134\u00d76=\"),\n    @(\"758\u00d76=\", \"445\u00d76=\"),\n    @(\"622\u00d78=\", \"595\u00d77=\"),\n    @(\"349\u00d75=\", \"539\u00d73=\"),\n    @(\"488\u00d79=\", \"396\u00d73=\"),\n    @(\"444\u00d76=\", \"308\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
